$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Shrub" was renamed to "CSS" and effectively rotated into group1, while
# "Grassland" moved from group1 to group2. The meandiff/lower/upper values
# correspondingly negate (and lower/upper swap) since the comparison
# direction between the two groups is reversed.
$ws.Range("A2").Value = "CSS"
$ws.Range("B2").Value = "Grassland"
$ws.Range("C2").Value = -2.0126
$ws.Range("E2").Value = -2.3017
$ws.Range("F2").Value = -1.7236
